# edit.ps1 — apply the "lead_v2 shielding thickness fix + running-thickness
# column" commit to the already-open workbook ($excel.ActiveWorkbook).
#
# Summary of the change (per the commit message "Updated the lead
# attenuation notebook, and fixed up the data for lead_v2 with the right
# shielding thicknesses"):
#   1. lead_v2 (sheet "lead_v2"): column B (lead thickness, mm) is
#      recomputed to be the *cumulative* plate thickness (running total)
#      instead of the single-plate thickness, for every data block.
#      The last two thickness tiers (which had no matching cumulative
#      measurement in plate_measurements) become 0.
#   2. plate_measurements: a new column F holding the running total
#      "D(this row) + D(previous row)" is added for the rows that feed the
#      lead_v2 cumulative thicknesses (rows 8-10), using a shared formula
#      for F9:F10.
#   3. Selection/active-cell bookkeeping on a couple of sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. lead_v2 — rewrite column B with the corrected cumulative thicknesses
# ---------------------------------------------------------------------
$leadV2 = $wb.Worksheets.Item("lead_v2")
$leadV2.Activate()

$leadV2.Range("B2:B43").Value   = 6.1724999999999994
$leadV2.Range("B44:B72").Value  = 12.614999999999998
$leadV2.Range("B73:B97").Value  = 26.442499999999999
$leadV2.Range("B98:B140").Value = 46.337500000000006
$leadV2.Range("B141:B199").Value = 0

$leadV2.Range("M115").Select()

# ---------------------------------------------------------------------
# 2. plate_measurements — add column F (running total of D)
# ---------------------------------------------------------------------
$plateMeas = $wb.Worksheets.Item("plate_measurements")
$plateMeas.Activate()

$plateMeas.Range("F8").Formula = "=D8+D7"
$plateMeas.Range("F9:F10").FormulaR1C1 = "=RC[-2]+R[-1]C[-2]"

$plateMeas.Range("F10").Select()

# ---------------------------------------------------------------------
# 3. lead_measures — selection only (kept on its previous active cell)
# ---------------------------------------------------------------------
$leadMeasures = $wb.Worksheets.Item("lead_measures")
$leadMeasures.Activate()
$leadMeasures.Range("C31").Select()

# Leave lead_v2 as the active sheet/tab, matching the target workbook
# (activeTab stays index 1 == lead_v2, and its sheetView keeps
# tabSelected="1").
$leadV2.Activate()

Write-Output "lead_v2 thicknesses + plate_measurements running totals updated"
